# Apply the checklist-update edit described by the commit:
# 1. Handled Same category in Axis category 1 and axis category 2 is not allowing drill down
# 2. Handled Date formatting
# 3. Handled cross filtering issue
#
# Concretely: two new checklist rows (#23 "Date Formatting", #24 "Drill down ")
# are appended to the "BVTs" sheet, each spanning 3 physical rows via merged
# cells (A/B/C vertically merged; D/E hold multiple step lines).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVTs")
$ws.Activate()

# ---------------------------------------------------------------
# Row block 88-90  -> item 23 "Date Formatting"
# ---------------------------------------------------------------
$ws.Range("A88").Value = 23
$ws.Range("B88").Value = "Date Formatting"
$ws.Range("C88").Value = "Check whether date formatting works"
$ws.Range("D88").Value = "1. Apply some format to date column from modelling tab`n2. Go to category filed in visual `n3. Drag that column in category field"
$ws.Range("E88").Value = "1.Date format is visilble in visual`n2.Format is change from modelling tab visual is updated"

$rngA88 = $ws.Range("A88:A90")
$rngA88.Merge()
$rngA88.VerticalAlignment = -4160   # xlVAlignTop
$rngA88.WrapText = $false

$rngB88 = $ws.Range("B88:B90")
$rngB88.Merge()
$rngB88.VerticalAlignment = -4160
$rngB88.WrapText = $true

$rngC88 = $ws.Range("C88:C90")
$rngC88.Merge()
$rngC88.Font.Name = "Calibri"

$rngD88 = $ws.Range("D88:D90")
$rngD88.Merge()
$rngD88.WrapText = $true

$rngE88 = $ws.Range("E88:E90")
$rngE88.Merge()
$rngE88.WrapText = $true

$ws.Rows(88).RowHeight = 14.25

# ---------------------------------------------------------------
# Row block 91-93  -> item 24 "Drill down "
# ---------------------------------------------------------------
$ws.Range("A91").Value = 24
$ws.Range("B91").Value = "Drill down "
$ws.Range("C91").Value = "Check whether drill down works when same Date field is dragged in Axis Category 1 and Axis Category 2"
$ws.Range("D91").Value = "1. Drag Date in Axis Category 1 (Date Hierarchy)"
$ws.Range("E91").Value = "1. Drill down works on Axis Category 1"
$ws.Range("D92").Value = "2. Drag Date in Axis Category 2"
$ws.Range("E92").Value = "2. Drill up works on Axis Category 1"
$ws.Range("D93").Value = "3. Drill Down on Axis Category 1"

$rngA91 = $ws.Range("A91:A93")
$rngA91.Merge()
$rngA91.VerticalAlignment = -4160
$rngA91.WrapText = $false

$rngB91 = $ws.Range("B91:B93")
$rngB91.Merge()
$rngB91.VerticalAlignment = -4160
$rngB91.WrapText = $true

$rngC91 = $ws.Range("C91:C93")
$rngC91.Merge()
$rngC91.Font.Name = "Calibri"

$rngD91 = $ws.Range("D91:D93")
$rngD91.VerticalAlignment = -4160
$rngD91.WrapText = $true

$rngE91 = $ws.Range("E91:E92")
$rngE91.VerticalAlignment = -4160
$rngE91.WrapText = $true

# Final selection / scroll position, matching the author's saved view.
try {
    $excel.ActiveWindow.ScrollRow = 82
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
}
$ws.Range("E92").Select()

Write-Host "Checklist rows 88-93 (items 23 and 24) added."
